$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of credential data
$ws.Range("A3").Value = "Driver"
$ws.Range("A4").Value = "Customer"
$ws.Range("B3").Value = "driver423"

$ws.Range("B4").Value = 4565678899
$ws.Range("B4").NumberFormat = "#,##0"

# Auto-fit column B like Excel would after entering data (bestFit width)
$ws.Columns.Item(2).AutoFit() | Out-Null

# Select B4 like the final state in the file
$ws.Range("B4").Select()
